# TC07_C3DC_phs000463_DiseasePhase-InitialDiag.xlsx
# Commit: "Updated remaining queries for C3DC"
#
# The DuckDB-style join queries stored in column B (rows 2-7) and C2
# referenced the raw PK/FK column name "id" on both sides of every JOIN
# (std.id / "study.id", prt.id / "participant.id"). They are updated to
# reference the fully-qualified column names study_id / participant_id
# that the source dataframes actually use.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells holding the SQL statements that need the join-column rename.
$queryCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $queryCells) {
    $rng = $ws.Range($addr)
    $txt = $rng.Value2
    if ($txt) {
        $txt = $txt.Replace('"participant.id"', '"participant.participant_id"')
        $txt = $txt.Replace('"study.id"', '"study.study_id"')
        $txt = $txt.Replace('prt.id', 'prt.participant_id')
        $txt = $txt.Replace('std.id', 'std.study_id')
        $rng.Value2 = $txt
    }
}

# Widen column C to fit the (now longer) StatQuery text and drop the
# stale "best fit" auto-size flag in favor of an explicit width.
$ws.Columns.Item(3).ColumnWidth = 69.25

# Refresh the view: scroll back to the top and move the selection/active
# cell from C7 to B2.
$excel.ActiveWindow.Zoom = 140
[void]$ws.Activate()
[void]$ws.Range("B2").Select()
